# Fix Senegal and Tanzania: remove the redundant "type" column (column G),
# which always held the constant value "character". Deleting the column
# shifts all subsequent columns (H -> G, I -> H, J -> I) to the left.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire column G (the "type" / "character" column).
$ws.Columns.Item(7).Delete()

# Match the recorded final selection from the edit.
$ws.Range("I7").Select()

$wb.Save()
